$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 243, pushing existing rows 243-351 down to 244-352
$ws.Rows.Item(243).Insert()

# Populate the new row 243 with the new data record
$ws.Cells.Item(243, 1).Value = 4
$ws.Cells.Item(243, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(243, 3).Value = "Los Lagos"
$ws.Cells.Item(243, 4).Value = 44636
$ws.Cells.Item(243, 5).Value = 10
$ws.Cells.Item(243, 6).Value = 100114001
$ws.Cells.Item(243, 7).Value = "Papa"
$ws.Cells.Item(243, 8).Value = "Patagonia"
$ws.Cells.Item(243, 9).Value = "1a (cosecha)"
$ws.Cells.Item(243, 10).Value = 250
$ws.Cells.Item(243, 11).Value = 7000
$ws.Cells.Item(243, 12).Value = 8000
$ws.Cells.Item(243, 13).Value = 7600
$ws.Cells.Item(243, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(243, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(243, 16).Value = 304
$ws.Cells.Item(243, 17).Value = 25
$ws.Cells.Item(243, 18).Value = "Hortaliza"
